# Adds a new "2022-Q3" sheet (fund holdings detail) and records the new
# quarter's summary row on the "总计" (total) sheet.

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    # Forces $value to be written as a text (string) cell, then re-applies
    # the "plain" style (copied from $styleSource) so no stray number-format
    # style is left behind on the cell.
    param($ws, $cellRef, $value, $styleSource)
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = $styleSource.Style
}

# ---------------------------------------------------------------------
# 1) 总计 (totals) sheet: insert a new row for 2022-Q3 at the top of the
#    data (row 2), pushing the existing quarters down by one row, and
#    renumber the sequential index column (A).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

# Restore the index-column style (bold/centered/bordered) on the new A2,
# copying it from the row right below (still carrying the original style).
$total.Range("A3").Copy($total.Range("A2"))

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 7
$total.Range("D2").Value = 0.93

# Renumber the 0-based sequence index for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
$total.Range("A8").Value = 6

# ---------------------------------------------------------------------
# 2) Add the new "2022-Q3" worksheet (fund holdings detail), positioned
#    right after "总计" and before "2022-Q2". Copying the "2022-Q2"
#    sheet gives us an identical layout/header/styling to start from.
# ---------------------------------------------------------------------
$sourceQ2 = $wb.Worksheets.Item("2022-Q2")
$sourceQ2.Copy($sourceQ2)

$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# Style reference cells (cells that already carry the correct "plain",
# un-styled look used by the data cells).
$plainIndex = $q3.Range("C2")

# Row 2
Set-TextCell $q3 "B2" "501011" $plainIndex
$q3.Range("C2").Value = "汇添富中证中药指数（LOF）A"
Set-TextCell $q3 "D2" "11.36" $plainIndex
Set-TextCell $q3 "E2" "94.73" $plainIndex
Set-TextCell $q3 "F2" "3.03" $plainIndex
Set-TextCell $q3 "G2" "0.3442" $plainIndex
$q3.Range("H2").Value = 9

# Row 3
Set-TextCell $q3 "B3" "501012" $plainIndex
$q3.Range("C3").Value = "汇添富中证中药指数（LOF）C"
Set-TextCell $q3 "D3" "6.42" $plainIndex
Set-TextCell $q3 "E3" "94.73" $plainIndex
Set-TextCell $q3 "F3" "3.03" $plainIndex
Set-TextCell $q3 "G3" "0.1945" $plainIndex
$q3.Range("H3").Value = 9

# Row 4
Set-TextCell $q3 "B4" "159647" $plainIndex
$q3.Range("C4").Value = "鹏华中证中药ETF"
Set-TextCell $q3 "D4" "6.16" $plainIndex
Set-TextCell $q3 "E4" "94.79" $plainIndex
Set-TextCell $q3 "F4" "2.99" $plainIndex
Set-TextCell $q3 "G4" "0.1842" $plainIndex
$q3.Range("H4").Value = 9

# Row 5 (new row, copy the A-column index style first)
$q3.Range("A4").Copy($q3.Range("A5"))
$q3.Range("A5").Value = 3
Set-TextCell $q3 "B5" "016950" $plainIndex
$q3.Range("C5").Value = "鹏华睿投灵活配置混合C"
Set-TextCell $q3 "D5" "4.12" $plainIndex
Set-TextCell $q3 "E5" "83.97" $plainIndex
Set-TextCell $q3 "F5" "1.83" $plainIndex
Set-TextCell $q3 "G5" "0.0754" $plainIndex
$q3.Range("H5").Value = 10

# Row 6 (new row)
$q3.Range("A4").Copy($q3.Range("A6"))
$q3.Range("A6").Value = 4
Set-TextCell $q3 "B6" "562390" $plainIndex
$q3.Range("C6").Value = "银华中证中药ETF"
Set-TextCell $q3 "D6" "2.34" $plainIndex
Set-TextCell $q3 "E6" "98.09" $plainIndex
Set-TextCell $q3 "F6" "3.14" $plainIndex
Set-TextCell $q3 "G6" "0.0735" $plainIndex
$q3.Range("H6").Value = 9

# Row 7 (new row)
$q3.Range("A4").Copy($q3.Range("A7"))
$q3.Range("A7").Value = 5
Set-TextCell $q3 "B7" "561510" $plainIndex
$q3.Range("C7").Value = "华泰柏瑞中证中药ETF"
Set-TextCell $q3 "D7" "2.02" $plainIndex
Set-TextCell $q3 "E7" "95.98" $plainIndex
Set-TextCell $q3 "F7" "3.02" $plainIndex
Set-TextCell $q3 "G7" "0.0610" $plainIndex
$q3.Range("H7").Value = 9

# Row 8 (new row) - note G8 is a genuine numeric 0, unlike the other
# (text-typed) G-column cells above.
$q3.Range("A4").Copy($q3.Range("A8"))
$q3.Range("A8").Value = 6
Set-TextCell $q3 "B8" "005434" $plainIndex
$q3.Range("C8").Value = "鹏华睿投灵活配置混合A"
Set-TextCell $q3 "D8" "0.00" $plainIndex
Set-TextCell $q3 "E8" "83.97" $plainIndex
Set-TextCell $q3 "F8" "1.83" $plainIndex
$q3.Range("G8").Value = 0
$q3.Range("H8").Value = 10

Write-Output "2022-Q3 sheet added and 总计 sheet updated"
